$wb = $excel.ActiveWorkbook

function Set-EventRow {
    param($ws, $row, $b, $c, $d, $e, $f, $g, $h, $i)
    $bCell = $ws.Cells.Item($row, 2)
    $bCell.NumberFormat = "@"
    $bCell.Value = $b
    $bCell.Style = "Normal"
    $ws.Cells.Item($row, 3).Value = $c
    $ws.Cells.Item($row, 4).Value = $d
    $eCell = $ws.Cells.Item($row, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $e
    $eCell.Style = "Normal"
    $ws.Cells.Item($row, 6).Value = $f
    $ws.Cells.Item($row, 7).Value = $g
    $ws.Cells.Item($row, 8).Value = $h
    $ws.Cells.Item($row, 9).Value = $i
}

$ws = $wb.Worksheets.Item("展览")

Set-EventRow $ws 2 '2024-07-28' '丽水·thp01～风摄少微' '大猷街 应星楼' '2024.07.28 10:00-07.28 18:00' 33 50 'https://show.bilibili.com/platform/detail.html?id=87134' '//i2.hdslb.com/bfs/openplatform/202407/WbSdFFLd1721636456044.jpeg'
Set-EventRow $ws 3 '2024-08-03' '丽水·樱卡动漫游戏嘉年华' '中东路848号(解放街交汇) 飞达国际大酒店' '2024.08.03 10:00-08.03 17:00' 273 50 'https://show.bilibili.com/platform/detail.html?id=87276' '//i0.hdslb.com/bfs/openplatform/202406/bVp0Zg1B1718172430380.jpeg'
Set-EventRow $ws 4 '2024-08-03' '丽水·逆光ZERO动漫游戏展' '丽阳街651号 丽水华侨君澜大饭店' '2024.08.03 10:00-08.03 17:00' 6 40 'https://show.bilibili.com/platform/detail.html?id=89687' '//i1.hdslb.com/bfs/openplatform/202407/NMYT1LRl1721639164353.jpeg'
Set-EventRow $ws 5 '2024-08-10' '丽水·CCAC动漫七夕（回馈展）' '中东路848号(解放街交汇) 飞达国际大酒店' '2024.08.10 09:00-08.10 17:00' 90 29.9 'https://show.bilibili.com/platform/detail.html?id=86567' '//i0.hdslb.com/bfs/openplatform/202405/tsOzbBRx1717015539538.png'
Set-EventRow $ws 6 '2024-08-17' '丽水·AEO纯白礼赞动漫嘉年华' '城北街1001号 爱依·时尚婚宴中心' '2024.08.17 09:00-08.17 18:00' 788 55 'https://show.bilibili.com/platform/detail.html?id=86779' '//i2.hdslb.com/bfs/openplatform/202406/MxJ3oNjt1717405405850.jpeg'
Set-EventRow $ws 7 '2024-08-24' '丽水·R动漫嘉年华' '中东路848号(解放街交汇) 飞达国际大酒店' '2024.08.24 09:30-08.24 17:00' 5 45 'https://show.bilibili.com/platform/detail.html?id=89651' '//i0.hdslb.com/bfs/openplatform/202407/7o5ALbAM1721383424201.jpeg'
Set-EventRow $ws 8 '2024-09-16' '丽水·LZ栗子动漫游戏嘉年华' '城北街798号 莱茵体育生活馆' '2024.09.16 09:30-09.16 17:00' 418 65 'https://show.bilibili.com/platform/detail.html?id=87480' '//i1.hdslb.com/bfs/openplatform/202406/bATqcZhH1719285865931.jpeg'

$ws.Rows.Item(9).Delete()

$ws = $wb.Worksheets.Item("全部类型")

Set-EventRow $ws 2 '2024-07-28' '丽水·thp01～风摄少微' '大猷街 应星楼' '2024.07.28 10:00-07.28 18:00' 33 50 'https://show.bilibili.com/platform/detail.html?id=87134' '//i2.hdslb.com/bfs/openplatform/202407/WbSdFFLd1721636456044.jpeg'
Set-EventRow $ws 3 '2024-08-03' '丽水·樱卡动漫游戏嘉年华' '中东路848号(解放街交汇) 飞达国际大酒店' '2024.08.03 10:00-08.03 17:00' 273 50 'https://show.bilibili.com/platform/detail.html?id=87276' '//i0.hdslb.com/bfs/openplatform/202406/bVp0Zg1B1718172430380.jpeg'
Set-EventRow $ws 4 '2024-08-03' '丽水·逆光ZERO动漫游戏展' '丽阳街651号 丽水华侨君澜大饭店' '2024.08.03 10:00-08.03 17:00' 6 40 'https://show.bilibili.com/platform/detail.html?id=89687' '//i1.hdslb.com/bfs/openplatform/202407/NMYT1LRl1721639164353.jpeg'
Set-EventRow $ws 5 '2024-08-10' '丽水·CCAC动漫七夕（回馈展）' '中东路848号(解放街交汇) 飞达国际大酒店' '2024.08.10 09:00-08.10 17:00' 90 29.9 'https://show.bilibili.com/platform/detail.html?id=86567' '//i0.hdslb.com/bfs/openplatform/202405/tsOzbBRx1717015539538.png'
Set-EventRow $ws 6 '2024-08-17' '丽水·AEO纯白礼赞动漫嘉年华' '城北街1001号 爱依·时尚婚宴中心' '2024.08.17 09:00-08.17 18:00' 788 55 'https://show.bilibili.com/platform/detail.html?id=86779' '//i2.hdslb.com/bfs/openplatform/202406/MxJ3oNjt1717405405850.jpeg'
Set-EventRow $ws 7 '2024-08-24' '丽水·R动漫嘉年华' '中东路848号(解放街交汇) 飞达国际大酒店' '2024.08.24 09:30-08.24 17:00' 5 45 'https://show.bilibili.com/platform/detail.html?id=89651' '//i0.hdslb.com/bfs/openplatform/202407/7o5ALbAM1721383424201.jpeg'
Set-EventRow $ws 8 '2024-09-16' '丽水·LZ栗子动漫游戏嘉年华' '城北街798号 莱茵体育生活馆' '2024.09.16 09:30-09.16 17:00' 418 65 'https://show.bilibili.com/platform/detail.html?id=87480' '//i1.hdslb.com/bfs/openplatform/202406/bATqcZhH1719285865931.jpeg'

$ws.Rows.Item(9).Delete()

